# "a bit more update" - Work Breakdown Structure
#
# Section "5. Closing" (rows 38-41) is reworked:
#   - old 5.3 "Document report on departments' progress"      -> removed
#   - old 5.4 "Document report on obtained supplies/services" -> replaced by
#       new 5.4 "Document Procurement Report" (P-mark moves from col B to col C)
#   - old 5.5 "Document report on messeges/mails"              -> removed
#   - old 5.6 "Document report on losses"                      -> replaced by
#       new 5.6 "Document Report on Losses" (P-mark stays, ends up in col E)
#
# Net effect: rows 40 & 41 disappear entirely and the sheet's used range
# shrinks from A1:E41 to A1:E39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38: "5.3 ..." -> "5.4 Document Procurement Report", P mark B38 -> C38
$ws.Range("A38").Value = "5.4 Document Procurement Report"
$ws.Range("B38").Copy($ws.Range("C38"))
$ws.Range("B38").Clear()

# --- Row 39: "5.4 ..." -> "5.6 Document Report on Losses", P mark C39 -> E39
$ws.Range("A39").Value = "5.6 Document Report on Losses"
$ws.Range("C39").Copy($ws.Range("E39"))
$ws.Range("C39").Clear()

# --- Old rows 40 ("5.5 ...") and 41 ("5.6 ...") are dropped completely
$ws.Rows("40:41").Delete()

# Reflect where the author ended up after the edit
$ws.Range("A39").Select()
